# Agregando margen de victoria y net rating en v2
# Insert a new bulleted list item ("Simple Rating System del equipo.")
# right after the "Días de descanso de los últimos 5 partidos." item
# inside the V2 list, matching the surrounding list's style/numbering.

$d = $word.ActiveDocument

# Locate the paragraph that ends the V2 bullet about rest days.
$findRange = $d.Content
$findRange.Find.Execute("Días de descanso de los últimos 5 partidos.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0) | Out-Null

$anchorPara = $findRange.Paragraphs(1)

# Create a new paragraph right after it; Word inherits the paragraph's
# style/numbering (pStyle "Prrafodelista", numPr ilvl 0 / numId 1) and
# run formatting (sz/szCs 24) automatically.
$anchorPara.Range.InsertParagraphAfter()
$newPara = $anchorPara.Next()

# Populate the new paragraph with the three runs from the target markup.
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:pPr>' + `
        '<w:pStyle w:val="Prrafodelista"/>' + `
        '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' + `
        '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>' + `
    '</w:pPr>' + `
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Simple Rating System</w:t></w:r>' + `
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> del equipo</w:t></w:r>' + `
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>.</w:t></w:r>' + `
    '</w:p>'

$newPara.Range.InsertXML($newParaXml)
